# Update column G ("K") values on Sheet1 per regenerated save_data
# (K = strikeouts column, values recalculated from updated source stats)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newK = @{
    2  = 6
    3  = 8
    4  = 3
    5  = 3
    6  = 9
    7  = 9
    8  = 5
    9  = 3
    10 = 5
    11 = 2
    12 = 4
    13 = 9
    14 = 13
    15 = 7
    16 = 6
    17 = 6
    18 = 7
    19 = 11
    20 = 5
    21 = 3
    22 = 8
    23 = 8
    24 = 6
    25 = 3
    26 = 2
    27 = 10
    28 = 5
    29 = 7
    30 = 7
    31 = 7
    32 = 4
    33 = 4
    34 = 5
    35 = 3
    36 = 5
    37 = 4
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
